# Add AT-400 to the Air Tractor aircraft data table.
#
# A new row is inserted between the existing "AT-301" (row 2) and
# "AT-401B" (row 3) entries, pushing the remaining rows down by one.
# The new row holds: Model="AT-400", Capacity=400, Wingspan=45,
# 70% boom=31.5, Ideal height=11.25, Working speed=135.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Air Tractor")

# Insert a new row above the current row 3 (AT-401B), shifting
# AT-401B..AT-802A down to rows 4..10.
$ws.Rows.Item(3).Insert()

# Populate the new row with the AT-400 data.
$ws.Cells.Item(3, 1).Value = "AT-400"
$ws.Cells.Item(3, 2).Value = 400
$ws.Cells.Item(3, 3).Value = 45
$ws.Cells.Item(3, 4).Value = 31.5
$ws.Cells.Item(3, 5).Value = 11.25
$ws.Cells.Item(3, 6).Value = 135

# Match the author's final UI state: the Air Tractor sheet is the
# active/selected tab, with H10 as the selected cell.
$ws.Activate() | Out-Null
$ws.Range("H10").Select() | Out-Null
